$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the empty spacer row (row 3), which shifts header/data rows up by one
$ws.Rows.Item(3).Delete()

# Row 1 height change
$ws.Rows.Item(1).RowHeight = 51.75

# Freeze panes below the new header row (row 3) and set selection
$ws.Rows.Item(4).Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D8").Select()

Write-Output "done"
